# Auto-generated edit script: updates crypto price/volume table
# (and swaps FraxShare/TheSandbox row order) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every touched cell to Text format first so values such as
# '1.001', '1.0000', '30.284.21' etc. are stored as literal text
# (matching the source inlineStr cells) instead of being
# auto-coerced to numbers/dates by Excel's smart entry.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.284.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.879.32'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.44'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4806'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2879'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06594'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.884.99'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.84'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07372'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.192'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6609'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.266.17'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.0000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007714'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.454'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.145.93'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '192.75'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.190'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.429'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.08'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.31'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.941'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.266'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09157'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.050'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05066'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7441'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.142'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.714'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01831'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.633'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9151'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.080'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.50'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.887'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4330'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.707'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1355'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +9.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.25'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -9.56%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.74%  '
